$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2217573221757322
$ws.Range("C2").Value = 0.5062761506276151
$ws.Range("J2").Value = 0.02092050209205021
$ws.Range("O2").Value = 0.004184100418410041
$ws.Range("P2").Value = 0.1506276150627615
$ws.Range("S2").Value = 0.09623430962343096

$ws.Range("B3").Value = 0.024
$ws.Range("C3").Value = 0.024
$ws.Range("J3").Value = 0.04
$ws.Range("P3").Value = 0.72
$ws.Range("S3").Value = 0.192

$ws.Range("J4").Value = 0.09523809523809523
$ws.Range("P4").Value = 0.5714285714285714
$ws.Range("S4").Value = 0.3333333333333333

$ws.Range("S5").Value = 1

$ws.Range("B6").Value = 0.05928853754940711
$ws.Range("D6").Value = 0.007905138339920948
$ws.Range("F6").Value = 0.1067193675889328
$ws.Range("J6").Value = 0.2134387351778656
$ws.Range("O6").Value = 0.0158102766798419
$ws.Range("Q6").Value = 0.1343873517786561
$ws.Range("R6").Value = 0.06719367588932806
$ws.Range("S6").Value = 0.3952569169960474

$ws.Range("B7").Value = 0.12
$ws.Range("D7").Value = 0.005
$ws.Range("F7").Value = 0.08500000000000001
$ws.Range("J7").Value = 0.135
$ws.Range("O7").Value = 0.01
$ws.Range("Q7").Value = 0.16
$ws.Range("R7").Value = 0.08
$ws.Range("S7").Value = 0.405

$ws.Range("B8").Value = 0.08126410835214447
$ws.Range("D8").Value = 0.009029345372460496
$ws.Range("F8").Value = 0.0654627539503386
$ws.Range("J8").Value = 0.126410835214447
$ws.Range("O8").Value = 0.006772009029345372
$ws.Range("Q8").Value = 0.1625282167042889
$ws.Range("R8").Value = 0.07223476297968397
$ws.Range("S8").Value = 0.4762979683972912

$ws.Range("B9").Value = 0.07878787878787878
$ws.Range("D9").Value = 0.006060606060606061
$ws.Range("E9").Value = 0.006060606060606061
$ws.Range("F9").Value = 0.06666666666666667
$ws.Range("J9").Value = 0.1272727272727273
$ws.Range("O9").Value = 0.01818181818181818
$ws.Range("Q9").Value = 0.1333333333333333
$ws.Range("R9").Value = 0.07272727272727272
$ws.Range("S9").Value = 0.4909090909090909

$ws.Range("B10").Value = 0.08409506398537477
$ws.Range("D10").Value = 0.03290676416819013
$ws.Range("F10").Value = 0.09140767824497258
$ws.Range("J10").Value = 0.1252285191956124
$ws.Range("O10").Value = 0.02010968921389397
$ws.Range("Q10").Value = 0.2056672760511883
$ws.Range("R10").Value = 0.07586837294332724
$ws.Range("S10").Value = 0.3647166361974406

$ws.Range("G11").Value = 0.1443661971830986
$ws.Range("J11").Value = 0.07394366197183098
$ws.Range("K11").Value = 0.1971830985915493
$ws.Range("L11").Value = 0.5598591549295775
$ws.Range("S11").Value = 0.02464788732394366

$ws.Range("G12").Value = 0.7639751552795031
$ws.Range("J12").Value = 0.1987577639751553
$ws.Range("K12").Value = 0.006211180124223602
$ws.Range("L12").Value = 0.0124223602484472
$ws.Range("S12").Value = 0.01863354037267081

$ws.Range("G13").Value = 0.7777777777777778
$ws.Range("J13").Value = 0.2037037037037037
$ws.Range("S13").Value = 0.01851851851851852

$ws.Range("S14").Value = 1

$ws.Range("F15").Value = 0.004651162790697674
$ws.Range("H15").Value = 0.213953488372093
$ws.Range("I15").Value = 0.06976744186046512
$ws.Range("J15").Value = 0.3488372093023256
$ws.Range("K15").Value = 0.07906976744186046
$ws.Range("M15").Value = 0.009302325581395349
$ws.Range("O15").Value = 0.06511627906976744
$ws.Range("S15").Value = 0.2093023255813954

$ws.Range("F16").Value = 0.03448275862068965
$ws.Range("H16").Value = 0.2413793103448276
$ws.Range("I16").Value = 0.04137931034482759
$ws.Range("J16").Value = 0.3310344827586207
$ws.Range("K16").Value = 0.1103448275862069
$ws.Range("M16").Value = 0.01379310344827586
$ws.Range("O16").Value = 0.06206896551724138
$ws.Range("S16").Value = 0.1655172413793103

$ws.Range("F17").Value = 0.01550387596899225
$ws.Range("H17").Value = 0.20671834625323
$ws.Range("I17").Value = 0.06976744186046512
$ws.Range("J17").Value = 0.3850129198966408
$ws.Range("K17").Value = 0.09819121447028424
$ws.Range("M17").Value = 0.01808785529715762
$ws.Range("O17").Value = 0.07235142118863049
$ws.Range("S17").Value = 0.1343669250645995

$ws.Range("F18").Value = 0.03164556962025317
$ws.Range("H18").Value = 0.1518987341772152
$ws.Range("I18").Value = 0.0949367088607595
$ws.Range("J18").Value = 0.3670886075949367
$ws.Range("K18").Value = 0.1075949367088608
$ws.Range("M18").Value = 0.04430379746835443
$ws.Range("N18").Value = 0.006329113924050633
$ws.Range("O18").Value = 0.0949367088607595
$ws.Range("S18").Value = 0.1012658227848101

$ws.Range("F19").Value = 0.01791530944625407
$ws.Range("H19").Value = 0.2174267100977199
$ws.Range("I19").Value = 0.08469055374592833
$ws.Range("J19").Value = 0.3330618892508143
$ws.Range("K19").Value = 0.1107491856677524
$ws.Range("M19").Value = 0.00244299674267101
$ws.Range("O19").Value = 0.06921824104234528
$ws.Range("S19").Value = 0.1327361563517915

